$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell H10 currently holds "type: blog`nwidth: 2`nheight: 1`nser: 168"
# Update the ser value from 168 to 175, keeping the rest of the text intact.
$ws.Range("H10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 175"
